$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6's shared string ("null") needs to be registered before row 1's
# ("page") so the rebuilt shared-strings table lands in the same order
# the workbook was saved with upstream ("a","null","page").
$ws.Range("A6").Value = "null"
$ws.Range("A6").HorizontalAlignment = -4152

$ws.Range("A1").Value = "page"

$ws.Range("A2").HorizontalAlignment = -4152

$ws.Range("A3").Value = -1
$ws.Range("A3").HorizontalAlignment = -4152

$ws.Range("A5").Value = 9999
$ws.Range("A5").HorizontalAlignment = -4152

$ws.Range("B3").Select()
